$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '67.452.13'
$ws.Range("E2").Value = '  -1.83%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.322.50'
$ws.Range("E3").Value = '  -2.35%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '578.38'
$ws.Range("E5").Value = '  -2.90%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '173.53'
$ws.Range("E6").Value = '  -7.41%  '
$ws.Range("E7").Value = '  +0.06%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.587'
$ws.Range("E8").Value = '  -2.62%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '3.320.38'
$ws.Range("E9").Value = '  -2.14%  '
$ws.Range("E10").Value = '  -5.60%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.574'
$ws.Range("E11").Value = '  -2.96%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '45.27'
$ws.Range("E12").Value = '  -5.14%  '
$ws.Range("E13").Value = '  -4.87%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '665.38'
$ws.Range("E14").Value = '  +3.41%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.863.37'
$ws.Range("E15").Value = '  -2.20%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '8.37'
$ws.Range("E16").Value = '  -3.28%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '67.688.55'
$ws.Range("E17").Value = '  -1.66%  '
$ws.Range("E18").Value = '  -1.29%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '3.327.22'
$ws.Range("E19").Value = '  -2.55%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '17.38'
$ws.Range("E20").Value = '  -4.12%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '10.88'
$ws.Range("E21").Value = '  -2.60%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.886'
$ws.Range("E22").Value = '  -3.31%  '
$ws.Range("E23").Value = '  +5.26%  '
$ws.Range("E24").Value = '  -6.67%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '98.29'
$ws.Range("E25").Value = '  -2.05%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '3.82'
$ws.Range("E26").Value = '  -7.03%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.66'
$ws.Range("E27").Value = '  -7.26%  '
$ws.Range("B28").Value = 'EthereumClassic'
$ws.Range("C28").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '33.77'
$ws.Range("E28").Value = '  +2.56%  '
$ws.Range("B29").Value = 'RenderToken'
$ws.Range("C29").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '9.23'
$ws.Range("E29").Value = '  -6.37%  '
$ws.Range("B30").Value = 'Filecoin'
$ws.Range("C30").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '8.38'
$ws.Range("E30").Value = '  -4.20%  '
$ws.Range("B31").Value = 'NEARProtocol'
$ws.Range("C31").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '7.32'
$ws.Range("E31").Value = '  +6.37%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '590.37'
$ws.Range("E32").Value = '  -4.07%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '10.88'
$ws.Range("E33").Value = '  -2.61%  '
$ws.Range("E34").Value = '  -2.60%  '
$ws.Range("E35").Value = '  +0.22%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '3.689.87'
$ws.Range("E36").Value = '  -8.54%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '56.70'
$ws.Range("E37").Value = '  -0.32%  '
$ws.Range("E38").Value = '  -15.33%  '
$ws.Range("E39").Value = '  -0.87%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '32.92'
$ws.Range("E40").Value = '  -2.74%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.61'
$ws.Range("E41").Value = '  -7.37%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '3.09'
$ws.Range("E42").Value = '  -7.37%  '
$ws.Range("E43").Value = '  -4.40%  '
$ws.Range("D44").Value = '0.0₃0657'
$ws.Range("E44").Value = '  -7.97%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '3.25'
$ws.Range("E45").Value = '  -5.15%  '
$ws.Range("E46").Value = '  -5.31%  '
$ws.Range("B47").Value = 'ThetaToken'
$ws.Range("C47").Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.58'
$ws.Range("E47").Value = '  -1.96%  '
$ws.Range("B48").Value = 'Stellar'
$ws.Range("C48").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.127'
$ws.Range("E48").Value = '  -2.32%  '
$ws.Range("E49").Value = '  +0.04%  '
$ws.Range("E50").Value = '  -4.03%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '126.85'
$ws.Range("E51").Value = '  -2.54%  '
